$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1257.55
$ws.Range("I53").Value = 98.333336
$ws.Range("J53").Value = 1754.3572
$ws.Range("K53").Value = 98.333336
$ws.Range("L53").Value = 1754.3572
$ws.Range("M53").Value = 538.666664
$ws.Range("N53").Value = -3028.3572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 5028.811
$ws.Range("I70").Value = 4293.2104
$ws.Range("J70").Value = 5805.278
$ws.Range("K70").Value = 12879.6312
$ws.Range("L70").Value = 17415.834
$ws.Range("M70").Value = -12609.6312
$ws.Range("N70").Value = -17955.834

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 5028.811
$ws.Range("I73").Value = 4293.2104
$ws.Range("J73").Value = 5805.278
$ws.Range("K73").Value = 12879.6312
$ws.Range("L73").Value = 17415.834
$ws.Range("M73").Value = -11943.6312
$ws.Range("N73").Value = -19287.834

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 4475
$ws.Range("I131").Value = 4633.3335
$ws.Range("J131").Value = 4000
$ws.Range("K131").Value = 13900.0005
$ws.Range("L131").Value = 12000
$ws.Range("M131").Value = -8860.000499999998
$ws.Range("N131").Value = -22080

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3299.4119
$ws.Range("I137").Value = 1995.8182
$ws.Range("J137").Value = 5689.3335
$ws.Range("K137").Value = 5987.4546
$ws.Range("L137").Value = 17068.0005
$ws.Range("M137").Value = -3437.4546
$ws.Range("N137").Value = -22168.0005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3182.2808
$ws.Range("I138").Value = 2542.2173
$ws.Range("J138").Value = 3615.2646
$ws.Range("K138").Value = 7626.651899999999
$ws.Range("L138").Value = 10845.7938
$ws.Range("M138").Value = -2486.651899999999
$ws.Range("N138").Value = -21125.7938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 16919.076
$ws.Range("I2").Value = 891
$ws.Range("K2").Value = 891
$ws.Range("M2").Value = -778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4489.982
$ws.Range("I32").Value = 4580.4185
$ws.Range("J32").Value = 4190.846
$ws.Range("K32").Value = 4580.4185
$ws.Range("L32").Value = 4190.846
$ws.Range("M32").Value = -4293.4185
$ws.Range("N32").Value = -4764.846

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 71432930
$ws.Range("I45").Value = 111113090
$ws.Range("J45").Value = 8631.200000000001
$ws.Range("K45").Value = 111113090
$ws.Range("L45").Value = 8631.200000000001
$ws.Range("M45").Value = -111112713
$ws.Range("N45").Value = -9385.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3966.818
$ws.Range("I61").Value = 3966.818
$ws.Range("K61").Value = 3966.818
$ws.Range("M61").Value = -3754.818

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 19611056
$ws.Range("I74").Value = 23811390
$ws.Range("K74").Value = 23811390
$ws.Range("M74").Value = -23810516

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 19611056
$ws.Range("I77").Value = 23811390
$ws.Range("K77").Value = 119056950
$ws.Range("M77").Value = -119052582

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2808.6667
$ws.Range("J88").Value = 2606.818
$ws.Range("L88").Value = 2606.818
$ws.Range("N88").Value = -3418.818

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2808.6667
$ws.Range("J91").Value = 2606.818
$ws.Range("L91").Value = 2606.818
$ws.Range("N91").Value = -5414.818

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1677.6666
$ws.Range("I97").Value = 1654.1428
$ws.Range("J97").Value = 1710.6
$ws.Range("K97").Value = 1654.1428
$ws.Range("L97").Value = 1710.6
$ws.Range("M97").Value = -1158.1428
$ws.Range("N97").Value = -2702.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 16919.076
$ws.Range("I116").Value = 891
$ws.Range("K116").Value = 891
$ws.Range("M116").Value = 1403

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3041.3262
$ws.Range("I132").Value = 860.34283
$ws.Range("K132").Value = 2581.02849
$ws.Range("M132").Value = -51.02849000000015

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3966.818
$ws.Range("I136").Value = 3966.818
$ws.Range("K136").Value = 11900.454
$ws.Range("M136").Value = -9350.454000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 16919.076
$ws.Range("I3").Value = 891
$ws.Range("K3").Value = 891
$ws.Range("M3").Value = -777

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6067.1665
$ws.Range("J134").Value = 11801.4
$ws.Range("L134").Value = 35404.2
$ws.Range("N134").Value = -40474.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28946.043
$ws.Range("I31").Value = 4954.423
$ws.Range("J31").Value = 58649.953
$ws.Range("K31").Value = 4954.423
$ws.Range("L31").Value = 58649.953
$ws.Range("M31").Value = -4659.423
$ws.Range("N31").Value = -59239.953

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 28946.043
$ws.Range("I34").Value = 4954.423
$ws.Range("J34").Value = 58649.953
$ws.Range("K34").Value = 4954.423
$ws.Range("L34").Value = 58649.953
$ws.Range("M34").Value = -4752.423
$ws.Range("N34").Value = -59053.953

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1094.875
$ws.Range("I107").Value = 914.86664
$ws.Range("J107").Value = 1394.8889
$ws.Range("K107").Value = 914.86664
$ws.Range("L107").Value = 1394.8889
$ws.Range("M107").Value = 1005.13336
$ws.Range("N107").Value = -5234.8889

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 517.2
$ws.Range("I23").Value = 46.5
$ws.Range("J23").Value = 589.61536
$ws.Range("K23").Value = 139.5
$ws.Range("L23").Value = 1768.84608
$ws.Range("M23").Value = 95.5
$ws.Range("N23").Value = -2238.84608

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 1750
$ws.Range("I126").Value = 1750
$ws.Range("K126").Value = 5250
$ws.Range("M126").Value = -310

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 9227.923000000001
$ws.Range("I113").Value = 9123.5
$ws.Range("J113").Value = 9395
$ws.Range("K113").Value = 9123.5
$ws.Range("L113").Value = 9395
$ws.Range("M113").Value = -6953.5
$ws.Range("N113").Value = -13735

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 7363.75
$ws.Range("I126").Value = 1798.6666
$ws.Range("J126").Value = 10702.8
$ws.Range("K126").Value = 5395.9998
$ws.Range("L126").Value = 32108.4
$ws.Range("M126").Value = -2925.9998
$ws.Range("N126").Value = -37048.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 63571.285
$ws.Range("I139").Value = 59999
$ws.Range("J139").Value = 64166.668
$ws.Range("K139").Value = 59999
$ws.Range("L139").Value = 64166.668
$ws.Range("M139").Value = -54859
$ws.Range("N139").Value = -74446.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6185.7646
$ws.Range("I40").Value = 5384.5625
$ws.Range("K40").Value = 5384.5625
$ws.Range("M40").Value = -5248.5625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5474.5835
$ws.Range("I46").Value = 2924.75
$ws.Range("J46").Value = 6749.5
$ws.Range("K46").Value = 2924.75
$ws.Range("L46").Value = 6749.5
$ws.Range("M46").Value = -2736.75
$ws.Range("N46").Value = -7125.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4566.95
$ws.Range("I136").Value = 1527.75
$ws.Range("K136").Value = 4583.25
$ws.Range("M136").Value = -2033.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4490.8066
$ws.Range("I126").Value = 4316.346
$ws.Range("K126").Value = 12949.038
$ws.Range("M126").Value = -10479.038
